# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# This updates the "K" column (column G) values on Sheet1 to reflect the
# recomputed strikeout (K) counts rather than the previous "Strike#" values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of row number -> new value for column G ("K")
$kValues = [ordered]@{
    2  = 1
    3  = 0
    4  = 0
    5  = 2
    6  = 1
    7  = 1
    8  = 1
    9  = 3
    10 = 1
    11 = 4
    12 = 2
    13 = 0
    14 = 3
    15 = 3
    16 = 1
    17 = 3
    18 = 5
    19 = 0
    20 = 3
    21 = 1
    22 = 0
    23 = 2
    24 = 0
    25 = 1
    26 = 1
    27 = 0
    28 = 3
    29 = 0
    30 = 3
    31 = 6
    32 = 3
    33 = 5
    34 = 4
    35 = 1
    36 = 3
    37 = 3
    38 = 2
    39 = 2
    40 = 2
    41 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
